$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 44875
$ws.Range("M2").Value = 50

# Row 3 updates
$ws.Range("D3").Value = 44874
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1600
$ws.Range("T3").Value = 10

# Row 4 updates
$ws.Range("D4").Value = 44855
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 5 kilos"
$ws.Range("S4").Value = 3000
$ws.Range("T4").Value = 5
